# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.423.08'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '1.918.13'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.63%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.25'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4068'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08212'
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.26'
$ws.Range("D12").Value = '1.906.79'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.065'
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.246'
$ws.Range("E14").Value = '  +2.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.76'
$ws.Range("E15").Value = '  +1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06876'
$ws.Range("E16").Value = '  +2.86%  '
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001039'
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.60'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").Value = '29.431.35'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.667'
$ws.Range("E22").Value = '  +2.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.73'
$ws.Range("E23").Value = '  -0.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.186'
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("D25").Value = '2.171.32'
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.664'
$ws.Range("E26").Value = '  +8.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.00'
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.01'
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.116'
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.91'
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("E31").Value = '  -1.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09623'
$ws.Range("E32").Value = '  +1.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.645'
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.546'
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.374'
$ws.Range("E35").Value = '  -1.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02282'
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06102'
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.181'
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.88'
$ws.Range("E39").Value = '  +6.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.071'
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5967'
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1847'
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.280'
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.403'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.46'
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07599'
$ws.Range("E46").Value = '  -2.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5591'
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.954'
$ws.Range("E48").Value = '  +1.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '118.76'
$ws.Range("E49").Value = '  +4.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.429'
$ws.Range("E50").Value = '  +3.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.26'
$ws.Range("E51").Value = '  -0.10%  '
